# 220330_2309 .gitignore 설정(DS_Store) / 앱 아이콘 제작
#
# 1) Refresh the cached "updated automatically" date placeholder text
#    (slide master + every slide layout) from 2021-10-04 to 2022. 3. 30.
# 2) On slide 1: duplicate the app-icon picture ("그림 3") into a second
#    picture ("그림 4") and move the title textbox down next to the new
#    icon (re-seating its shape id from 6 -> 7 along the way).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder refresh across the slide master and all layouts
# ---------------------------------------------------------------------
$oldDate = "2021-10-04"
$newDate = "2022. 3. 30."

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1: add the second app-icon picture + move the title
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$icon = $s.Shapes.Item(1)   # "그림 3" picture (id 4)
$title = $s.Shapes.Item(2)  # "제목 1" textbox (id 6)

# Consume the two lowest free shape ids (2,3) with throw-away textboxes so
# that the ids handed out below line up with the ones PowerPoint used
# (icon duplicate -> 5, title re-seat -> 7) instead of the lowest-unused
# ids (2 / 3).
$spacer1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$spacer2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$spacer1.Delete()
$spacer2.Delete()

# Duplicate the icon picture -> new shape picks up id 5 ("그림 4")
$iconDup = $icon.Duplicate().Item(1)
$iconDup.Name = "그림 4"
$iconDup.Left = 5548790 / 12700
$iconDup.Top = 5210424 / 12700
$iconDup.Width = 857250 / 12700
$iconDup.Height = 857250 / 12700

# Re-seat the title: duplicate it (keeps text/formatting, new id 7),
# move the duplicate to its new position, then drop the original (id 6).
$titleDup = $title.Duplicate().Item(1)
$titleDup.Left = 4602705 / 12700
$titleDup.Top = 5349770 / 12700
$titleDup.Width = 2749420 / 12700
$titleDup.Height = 717904 / 12700
$title.Delete()
